$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting D:K to E:L
$ws.Columns("D:D").Insert()

# Copy number formats from column E (old column D, now shifted) into new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

# Populate new column D with FY2018 data
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 518900
$ws.Range("D9").Value = 330700
$ws.Range("D10").Value = 188200
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 448400
$ws.Range("D18").Value = 70500
$ws.Range("D20").Value = 2900
$ws.Range("D21").Value = 92100
$ws.Range("D22").Value = 200
$ws.Range("D23").Value = 73200
$ws.Range("D24").Value = 16400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 56800
$ws.Range("D27").Value = 56900
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -2900
$ws.Range("D33").Value = 56900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 56900
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 110900
$ws.Range("D42").Value = 75100
$ws.Range("D43").Value = 52600
$ws.Range("D44").Value = 54500
$ws.Range("D45").Value = 10900
$ws.Range("D46").Value = 304000
$ws.Range("D47").Value = 170400
$ws.Range("D48").Value = 186100
$ws.Range("D49").Value = 248300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 38500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 947400
$ws.Range("D57").Value = 11800
$ws.Range("D58").Value = 400
$ws.Range("D59").Value = 49200
$ws.Range("D60").Value = 61400
$ws.Range("D61").Value = 7500
$ws.Range("D62").Value = 128000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 196700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 33800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 750600
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 56900
$ws.Range("D83").Value = 18700
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 100900
$ws.Range("D91").Value = -27600
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -44500
$ws.Range("D96").Value = -23000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -42400
$ws.Range("D101").Value = 500
$ws.Range("D102").Value = 14600
